$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "37.817.92"
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = "  +0.01%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.089.16"
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = "  +0.03%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "233.76"
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = "  -0.58%  "
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = "  -0.08%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "58.61"
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = "  -0.36%  "
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = "  +0.03%  "
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = "  +0.64%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.0784"
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = "  -0.82%  "
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = "  +2.74%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "15.16"
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = "  +2.53%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "2.397.73"
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = "  +0.07%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "21.40"
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = "  +0.67%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.781"
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = "  +1.32%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "5.38"
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = "  +0.98%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "2.095.61"
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = "  +0.42%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "37.786.14"
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = "  +0.12%  "
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = "  -0.84%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "71.32"
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = "  -0.04%  "
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = "  +0.17%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "230.40"
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = "  +0.72%  "
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = "  -0.10%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "2.39"
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = "  -0.62%  "
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = "  -1.32%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "9.84"
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = "  +9.23%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "172.09"
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = "  +1.10%  "
$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = "  -1.93%  "
$ws.Range("E29").NumberFormat = "@"
$ws.Range("E29").Value = "  +0.06%  "
$ws.Range("E30").NumberFormat = "@"
$ws.Range("E30").Value = "  -0.70%  "
$ws.Range("E31").NumberFormat = "@"
$ws.Range("E31").Value = "  +0.87%  "
$ws.Range("E32").NumberFormat = "@"
$ws.Range("E32").Value = "  +0.56%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.0633"
$ws.Range("E33").NumberFormat = "@"
$ws.Range("E33").Value = "  +0.18%  "
$ws.Range("E34").NumberFormat = "@"
$ws.Range("E34").Value = "  -0.92%  "
$ws.Range("E35").NumberFormat = "@"
$ws.Range("E35").Value = "  -0.33%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "1.83"
$ws.Range("E37").NumberFormat = "@"
$ws.Range("E37").Value = "  -1.76%  "
$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = "  +0.15%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "5.37"
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = "  -0.31%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.0237"
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = "  +9.98%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "102.70"
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = "  +3.84%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.0974"
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = "  -2.73%  "
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = "  -0.31%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "16.82"
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = "  +4.83%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "1.455.63"
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = "  -0.70%  "
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = "  -0.78%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "4.20"
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = "  -4.25%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "1.06"
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = "  -1.30%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "7.34"
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = "  -0.89%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "2.98"
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = "  -2.11%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "2.281.59"
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = "  +0.06%  "
